# Reproduces the "QuizAnswer" alt-text tagging pass + the slide-18
# (sldId 516) grouping of the TextBox/Picture answer callout, as
# described by the commit's OOXML diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Tag every quiz "Smiley Face" answer-reveal shape with
#    descr="QuizAnswer" (Shape.AlternativeText), on the slides the
#    diff touches: sldId 497, 498, 505, 506, 508, 516.
# ---------------------------------------------------------------

# Slide 5 (sldId 497): shape id=5 "Smiley Face 4"
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(3).AlternativeText = "QuizAnswer"

# Slide 6 (sldId 498): shape id=5 "Smiley Face 4"
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(3).AlternativeText = "QuizAnswer"

# Slide 11 (sldId 505): shape id=32 "Smiley Face 31"
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item(22).AlternativeText = "QuizAnswer"

# Slide 15 (sldId 506): shape id=4 "Smiley Face 3"
$s15 = $p.Slides.Item(15)
$s15.Shapes.Item(3).AlternativeText = "QuizAnswer"

# Slide 17 (sldId 508): shape id=6 "Smiley Face 5"
$s17 = $p.Slides.Item(17)
$s17.Shapes.Item(4).AlternativeText = "QuizAnswer"

# ---------------------------------------------------------------
# 2) Slide 18 (sldId 516): tag "Smiley Face 5" + "TextBox 6", group
#    "TextBox 6" and "Picture 7" together, nudge the picture down,
#    and repoint / prune the click animations to match.
# ---------------------------------------------------------------

$s18 = $p.Slides.Item(18)

# shape id=6 "Smiley Face 5"
$s18.Shapes.Item(4).AlternativeText = "QuizAnswer"

# Locate "TextBox 6" (id=7) and "Picture 7" (id=8) by name so the
# script doesn't depend on a fragile fixed index order.
$txtIdx = 0
$picIdx = 0
for ($i = 1; $i -le $s18.Shapes.Count; $i++) {
    $nm = $s18.Shapes.Item($i).Name
    if ($nm -eq "TextBox 6") { $txtIdx = $i }
    if ($nm -eq "Picture 7") { $picIdx = $i }
}

$s18.Shapes.Item($txtIdx).AlternativeText = "QuizAnswer"

# Group the textbox + picture into a new "Group 2" shape.
$range = $s18.Shapes.Range(@($txtIdx, $picIdx))
$grp = $range.Group()
$grp.Name = "Group 2"

# Inside the group, item 1 is the textbox, item 2 is the picture
# (same relative order they were added to the Range in).
$pic = $grp.GroupItems.Item(2)
# Move the picture down from y=5967789 EMU to y=6027261 EMU.
$pic.Top = 6027261 / 914400 * 72

# Fix up the click-animation sequence: the picture's entrance effect
# now targets the new group, and the textbox's separate "with
# previous" entrance effect (which is now redundant, since it
# animates together with the picture inside the group) is removed.
$seq = $s18.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $eff = $seq.Item($i)
    $nm = $eff.Shape.Name
    if ($nm -eq "TextBox 6" -or $nm -eq "Picture 7") {
        $eff.Delete()
    }
}
$null = $seq.AddEffect($grp, 10)
